$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.6823026622293469
$ws.Cells.Item(2, 3).Value = 0.2211368881851605
$ws.Cells.Item(2, 4).Value = 0.5694247468477442
$ws.Cells.Item(2, 5).Value = 0.2064649584887519
$ws.Cells.Item(2, 7).Value = 1.426853873779166
$ws.Cells.Item(2, 8).Value = 1.340461861374479
$ws.Cells.Item(2, 10).Value = 0.09579734044755295
$ws.Cells.Item(2, 11).Value = 0.2972150827219195
$ws.Cells.Item(2, 12).Value = 0.3408298847538802
$ws.Cells.Item(2, 13).Value = 0.2303698922331421
$ws.Cells.Item(2, 15).Value = 5.662285315498082

# Row 3
$ws.Cells.Item(3, 2).Value = 0.6498435481162517
$ws.Cells.Item(3, 3).Value = 0.2213868195474689
$ws.Cells.Item(3, 4).Value = 0.5682856462176744
$ws.Cells.Item(3, 5).Value = 0.2072449609768547
$ws.Cells.Item(3, 7).Value = 1.434624654355545
$ws.Cells.Item(3, 8).Value = 1.348364743891352
$ws.Cells.Item(3, 10).Value = 0.09611060598067844
$ws.Cells.Item(3, 11).Value = 0.2679151331702485
$ws.Cells.Item(3, 12).Value = 0.3389894649816796
$ws.Cells.Item(3, 13).Value = 0.2244617232952599
$ws.Cells.Item(3, 15).Value = 5.695007427484541

# Row 4
$ws.Cells.Item(4, 2).Value = 0.6301418147683933
$ws.Cells.Item(4, 3).Value = 0.2215630116707032
$ws.Cells.Item(4, 4).Value = 0.567810077904312
$ws.Cells.Item(4, 5).Value = 0.2077861852938501
$ws.Cells.Item(4, 7).Value = 1.440031430009348
$ws.Cells.Item(4, 8).Value = 1.353658743524107
$ws.Cells.Item(4, 10).Value = 0.09631315480749825
$ws.Cells.Item(4, 11).Value = 0.2499550536106199
$ws.Cells.Item(4, 12).Value = 0.3379789310784957
$ws.Cells.Item(4, 13).Value = 0.2209194574215019
$ws.Cells.Item(4, 15).Value = 5.717360123319764

# Row 5
$ws.Cells.Item(5, 2).Value = 0.6221712633665391
$ws.Cells.Item(5, 3).Value = 0.2216405523672762
$ws.Cells.Item(5, 4).Value = 0.5676727104306138
$ws.Cells.Item(5, 5).Value = 0.2080224377558419
$ws.Cells.Item(5, 7).Value = 1.442394639127677
$ws.Cells.Item(5, 8).Value = 1.355927279958877
$ws.Cells.Item(5, 10).Value = 0.09639826647124705
$ws.Cells.Item(5, 11).Value = 0.2426442527462882
$ws.Cells.Item(5, 12).Value = 0.3375972853272202
$ws.Cells.Item(5, 13).Value = 0.2194975707905265
$ws.Cells.Item(5, 15).Value = 5.72703814428283

# Row 6
$ws.Cells.Item(6, 2).Value = 0.6208512857650987
$ws.Cells.Item(6, 3).Value = 0.2216537754049313
$ws.Cells.Item(6, 4).Value = 0.5676533131360628
$ws.Cells.Item(6, 5).Value = 0.2080626164703183
$ws.Cells.Item(6, 7).Value = 1.442796709146783
$ws.Cells.Item(6, 8).Value = 1.35631068849812
$ws.Cells.Item(6, 10).Value = 0.09641255467589405
$ws.Cells.Item(6, 11).Value = 0.2414308018841211
$ws.Cells.Item(6, 12).Value = 0.3375357378086292
$ws.Cells.Item(6, 13).Value = 0.2192627774669873
$ws.Cells.Item(6, 15).Value = 5.728679560860371

# Row 7
$ws.Cells.Item(7, 2).Value = 0.6300340851473436
$ws.Cells.Item(7, 3).Value = 0.2215640341335252
$ws.Cells.Item(7, 4).Value = 0.5678079966432534
$ws.Cells.Item(7, 5).Value = 0.2077893078674151
$ws.Cells.Item(7, 7).Value = 1.440062653488312
$ws.Cells.Item(7, 8).Value = 1.353688887409731
$ws.Cells.Item(7, 10).Value = 0.0963142922339153
$ws.Cells.Item(7, 11).Value = 0.2498564242491597
$ws.Cells.Item(7, 12).Value = 0.3379736618301692
$ws.Cells.Item(7, 13).Value = 0.2209001936448196
$ws.Cells.Item(7, 15).Value = 5.717488339477086

# Row 8
$ws.Cells.Item(8, 2).Value = 0.6710637606825856
$ws.Cells.Item(8, 3).Value = 0.2212183608131575
$ws.Cells.Item(8, 4).Value = 0.5689856103862354
$ws.Cells.Item(8, 5).Value = 0.2067209929795215
$ws.Cells.Item(8, 7).Value = 1.429401406783143
$ws.Cells.Item(8, 8).Value = 1.343095208637337
$ws.Cells.Item(8, 10).Value = 0.09590324041457876
$ws.Cells.Item(8, 11).Value = 0.2871064696271475
$ws.Cells.Item(8, 12).Value = 0.3401705711096099
$ws.Cells.Item(8, 13).Value = 0.22831512242362
$ws.Cells.Item(8, 15).Value = 5.673098901218907

# Row 9
$ws.Cells.Item(9, 2).Value = 0.7533098098158462
$ws.Cells.Item(9, 3).Value = 0.2207198278303935
$ws.Cells.Item(9, 4).Value = 0.5730653758840418
$ws.Cells.Item(9, 5).Value = 0.2051189324704854
$ws.Cells.Item(9, 7).Value = 1.41353291950702
$ws.Cells.Item(9, 8).Value = 1.325818625838124
$ws.Cells.Item(9, 10).Value = 0.09517785035156168
$ws.Cells.Item(9, 11).Value = 0.3603759037809198
$ws.Cells.Item(9, 12).Value = 0.3454225712784478
$ws.Cells.Item(9, 13).Value = 0.2435276144065099
$ws.Cells.Item(9, 15).Value = 5.603971310345713

# Row 10
$ws.Cells.Item(10, 2).Value = 0.8147991697789507
$ws.Cells.Item(10, 3).Value = 0.220461517300734
$ws.Cells.Item(10, 4).Value = 0.5771354923438281
$ws.Cells.Item(10, 5).Value = 0.204240551892358
$ws.Cells.Item(10, 7).Value = 1.404941079910557
$ws.Cells.Item(10, 8).Value = 1.315249419026244
$ws.Cells.Item(10, 10).Value = 0.09469370724216919
$ws.Cells.Item(10, 11).Value = 0.4143243363541274
$ws.Cells.Item(10, 12).Value = 0.3498516598138082
$ws.Cells.Item(10, 13).Value = 0.2551076805717614
$ws.Cells.Item(10, 15).Value = 5.564080878544388

# Row 11
$ws.Cells.Item(11, 2).Value = 0.8429979970921977
$ws.Cells.Item(11, 3).Value = 0.2203671656957411
$ws.Cells.Item(11, 4).Value = 0.5792186833718063
$ws.Cells.Item(11, 5).Value = 0.2039054218025207
$ws.Cells.Item(11, 7).Value = 1.401697422324105
$ws.Cells.Item(11, 8).Value = 1.31090066179334
$ws.Cells.Item(11, 10).Value = 0.09448397284552357
$ws.Cells.Item(11, 11).Value = 0.4388888712294943
$ws.Cells.Item(11, 12).Value = 0.3519894302250179
$ws.Cells.Item(11, 13).Value = 0.2604621455223111
$ws.Cells.Item(11, 15).Value = 5.54829451313816

# Row 12
$ws.Cells.Item(12, 2).Value = 0.853708227237064
$ws.Cells.Item(12, 3).Value = 0.2203347435206524
$ws.Cells.Item(12, 4).Value = 0.5800407178393527
$ws.Cells.Item(12, 5).Value = 0.2037877523647325
$ws.Cells.Item(12, 7).Value = 1.400564647673619
$ws.Cells.Item(12, 8).Value = 1.309319799142344
$ws.Cells.Item(12, 10).Value = 0.09440605644055067
$ws.Cells.Item(12, 11).Value = 0.4481937329580603
$ws.Cells.Item(12, 12).Value = 0.35281653084526
$ws.Cells.Item(12, 13).Value = 0.2625020695070148
$ws.Cells.Item(12, 15).Value = 5.542655529317585

# Row 13
$ws.Cells.Item(13, 2).Value = 0.8514001778162879
$ws.Cells.Item(13, 3).Value = 0.2203415795101762
$ws.Cells.Item(13, 4).Value = 0.5798622047500714
$ws.Cells.Item(13, 5).Value = 0.2038126842515275
$ws.Cells.Item(13, 7).Value = 1.400804363461802
$ws.Cells.Item(13, 8).Value = 1.309657336397819
$ws.Cells.Item(13, 10).Value = 0.09442277025762058
$ws.Cells.Item(13, 11).Value = 0.4461896479807308
$ws.Cells.Item(13, 12).Value = 0.3526376200584593
$ws.Cells.Item(13, 13).Value = 0.2620621905840892
$ws.Cells.Item(13, 15).Value = 5.543854915849948

# Row 14
$ws.Cells.Item(14, 2).Value = 0.8438784979106515
$ws.Cells.Item(14, 3).Value = 0.2203644321624054
$ws.Cells.Item(14, 4).Value = 0.5792856485600169
$ws.Cells.Item(14, 5).Value = 0.2038955561184288
$ws.Cells.Item(14, 7).Value = 1.401602314155767
$ws.Cells.Item(14, 8).Value = 1.310769282869444
$ws.Cells.Item(14, 10).Value = 0.09447753247338309
$ws.Cells.Item(14, 11).Value = 0.4396543345485782
$ws.Cells.Item(14, 12).Value = 0.3520571247095603
$ws.Cells.Item(14, 13).Value = 0.2606297255785819
$ws.Cells.Item(14, 15).Value = 5.547823799331383

# Row 15
$ws.Cells.Item(15, 2).Value = 0.8392753914554589
$ws.Cells.Item(15, 3).Value = 0.2203788600289514
$ws.Cells.Item(15, 4).Value = 0.5789368075061532
$ws.Cells.Item(15, 5).Value = 0.2039475195566069
$ws.Cells.Item(15, 7).Value = 1.402103520453892
$ws.Cells.Item(15, 8).Value = 1.311458963050327
$ws.Cells.Item(15, 10).Value = 0.09451127181657348
$ws.Cells.Item(15, 11).Value = 0.4356516163324784
$ws.Cells.Item(15, 12).Value = 0.3517038399297832
$ws.Cells.Item(15, 13).Value = 0.259753897230496
$ws.Cells.Item(15, 15).Value = 5.550298987310043

# Row 16
$ws.Cells.Item(16, 2).Value = 0.812960818450648
$ws.Cells.Item(16, 3).Value = 0.2204681473732961
$ws.Cells.Item(16, 4).Value = 0.5770040008896729
$ws.Cells.Item(16, 5).Value = 0.2042637476949025
$ws.Cells.Item(16, 7).Value = 1.405166432217442
$ws.Cells.Item(16, 8).Value = 1.315542851184389
$ws.Cells.Item(16, 10).Value = 0.09470762482205064
$ws.Cells.Item(16, 11).Value = 0.4127194070709663
$ws.Cells.Item(16, 12).Value = 0.3497144164784771
$ws.Cells.Item(16, 13).Value = 0.2547594838660743
$ws.Cells.Item(16, 15).Value = 5.565160010034504

# Row 17
$ws.Cells.Item(17, 2).Value = 0.7968753260268215
$ws.Cells.Item(17, 3).Value = 0.2205288364413107
$ws.Cells.Item(17, 4).Value = 0.5758775249557573
$ws.Cells.Item(17, 5).Value = 0.2044742286230274
$ws.Cells.Item(17, 7).Value = 1.407215653362513
$ws.Cells.Item(17, 8).Value = 1.318165718631505
$ws.Cells.Item(17, 10).Value = 0.0948307678934448
$ws.Cells.Item(17, 11).Value = 0.3986567915748083
$ws.Cells.Item(17, 12).Value = 0.3485253849869423
$ws.Cells.Item(17, 13).Value = 0.2517176528571454
$ws.Cells.Item(17, 15).Value = 5.574880936178147

# Row 18
$ws.Cells.Item(18, 2).Value = 0.7876447849893111
$ws.Cells.Item(18, 3).Value = 0.2205659240145366
$ws.Cells.Item(18, 4).Value = 0.5752514204357198
$ws.Cells.Item(18, 5).Value = 0.2046013598227745
$ws.Cells.Item(18, 7).Value = 1.408456890809532
$ws.Cells.Item(18, 8).Value = 1.319717554396973
$ws.Cells.Item(18, 10).Value = 0.09490258558843445
$ws.Cells.Item(18, 11).Value = 0.3905705544260911
$ws.Cells.Item(18, 12).Value = 0.3478530651229335
$ws.Cells.Item(18, 13).Value = 0.2499762357568613
$ws.Cells.Item(18, 15).Value = 5.58069431287646

# Row 19
$ws.Cells.Item(19, 2).Value = 0.784523180645806
$ws.Cells.Item(19, 3).Value = 0.220578856493745
$ws.Cells.Item(19, 4).Value = 0.5750431831325074
$ws.Cells.Item(19, 5).Value = 0.2046454474085735
$ws.Cells.Item(19, 7).Value = 1.408887902767177
$ws.Cells.Item(19, 8).Value = 1.32025040806468
$ws.Cells.Item(19, 10).Value = 0.09492707185398785
$ws.Cells.Item(19, 11).Value = 0.3878330884326147
$ws.Cells.Item(19, 12).Value = 0.3476274212565897
$ws.Cells.Item(19, 13).Value = 0.2493880286484682
$ws.Cells.Item(19, 15).Value = 5.582700790058169

# Row 20
$ws.Cells.Item(20, 2).Value = 0.798585443485706
$ws.Cells.Item(20, 3).Value = 0.2205221504447721
$ws.Cells.Item(20, 4).Value = 0.5759951835172643
$ws.Cells.Item(20, 5).Value = 0.2044511947624379
$ws.Cells.Item(20, 7).Value = 1.406991034029474
$ws.Cells.Item(20, 8).Value = 1.317882036450925
$ws.Cells.Item(20, 10).Value = 0.09481755677697867
$ws.Cells.Item(20, 11).Value = 0.4001535563170648
$ws.Cells.Item(20, 12).Value = 0.3486507617538734
$ws.Cells.Item(20, 13).Value = 0.2520406171686815
$ws.Cells.Item(20, 15).Value = 5.573823137421385

# Row 21
$ws.Cells.Item(21, 2).Value = 0.8460869361663583
$ws.Cells.Item(21, 3).Value = 0.2203576302223667
$ws.Cells.Item(21, 4).Value = 0.5794540978322971
$ws.Cells.Item(21, 5).Value = 0.2038709641956977
$ws.Cells.Item(21, 7).Value = 1.401365344658217
$ws.Cells.Item(21, 8).Value = 1.310440889205083
$ws.Cells.Item(21, 10).Value = 0.09446140666553138
$ws.Cells.Item(21, 11).Value = 0.441573842514515
$ws.Cells.Item(21, 12).Value = 0.3522271541853144
$ws.Cells.Item(21, 13).Value = 0.2610501425618352
$ws.Cells.Item(21, 15).Value = 5.546648845626891

# Row 22
$ws.Cells.Item(22, 2).Value = 0.8773177921203512
$ws.Cells.Item(22, 3).Value = 0.2202693704348704
$ws.Cells.Item(22, 4).Value = 0.581908004307877
$ws.Cells.Item(22, 5).Value = 0.2035455764030338
$ws.Cells.Item(22, 7).Value = 1.398245409021911
$ws.Cells.Item(22, 8).Value = 1.305961831191937
$ws.Cells.Item(22, 10).Value = 0.09423741461400548
$ws.Cells.Item(22, 11).Value = 0.4686605064436549
$ws.Cells.Item(22, 12).Value = 0.3546669146140431
$ws.Cells.Item(22, 13).Value = 0.2670100360822403
$ws.Cells.Item(22, 15).Value = 5.530864470437308

# Row 23
$ws.Cells.Item(23, 2).Value = 0.8606325181299894
$ws.Cells.Item(23, 3).Value = 0.2203147212320857
$ws.Cells.Item(23, 4).Value = 0.5805806667500235
$ws.Cells.Item(23, 5).Value = 0.2037143269987354
$ws.Cells.Item(23, 7).Value = 1.399859655791815
$ws.Cells.Item(23, 8).Value = 1.308317277090097
$ws.Cells.Item(23, 10).Value = 0.09435616236534816
$ws.Cells.Item(23, 11).Value = 0.4542025411135739
$ws.Cells.Item(23, 12).Value = 0.353355436147524
$ws.Cells.Item(23, 13).Value = 0.2638226236598982
$ws.Cells.Item(23, 15).Value = 5.539108251301741

# Row 24
$ws.Cells.Item(24, 2).Value = 0.797812245683275
$ws.Cells.Item(24, 3).Value = 0.2205251663393
$ws.Cells.Item(24, 4).Value = 0.5759419230354581
$ws.Cells.Item(24, 5).Value = 0.2044615892963932
$ws.Cells.Item(24, 7).Value = 1.407092387831185
$ws.Cells.Item(24, 8).Value = 1.318010152363939
$ws.Cells.Item(24, 10).Value = 0.09482352634178737
$ws.Cells.Item(24, 11).Value = 0.3994768734621914
$ws.Cells.Item(24, 12).Value = 0.3485940438043116
$ws.Cells.Item(24, 13).Value = 0.2518945820289602
$ws.Cells.Item(24, 15).Value = 5.574300668320035

# Row 25
$ws.Cells.Item(25, 2).Value = 0.7308714810504
$ws.Cells.Item(25, 3).Value = 0.2208356375225975
$ws.Cells.Item(25, 4).Value = 0.5717728341252979
$ws.Cells.Item(25, 5).Value = 0.205499753598188
$ws.Cells.Item(25, 7).Value = 1.417286876715508
$ws.Cells.Item(25, 8).Value = 1.330118821657763
$ws.Cells.Item(25, 10).Value = 0.09536549003998873
$ws.Cells.Item(25, 11).Value = 0.3405326671320097
$ws.Cells.Item(25, 12).Value = 0.3439012324463562
$ws.Cells.Item(25, 13).Value = 0.2393409414440946
$ws.Cells.Item(25, 15).Value = 5.620756484457587
